$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(12, 8).Value = 2612.75  # H12: 2117 -> 2612.75
$ws.Cells.Item(12, 9).Value = 2612.75  # I12: 2117 -> 2612.75
$ws.Cells.Item(12, 11).Value = 2612.75  # K12: 2117 -> 2612.75
$ws.Cells.Item(12, 13).Value = -2442.75  # M12: -1947 -> -2442.75
$ws.Cells.Item(86, 8).Value = 5148.9443  # H86: 5042.684 -> 5148.9443
$ws.Cells.Item(86, 9).Value = 4697.125  # I86: 4523 -> 4697.125
$ws.Cells.Item(86, 11).Value = 4697.125  # K86: 4523 -> 4697.125
$ws.Cells.Item(86, 13).Value = -3574.125  # M86: -3400 -> -3574.125
$ws.Cells.Item(88, 8).Value = 5983  # H88: 6099.6665 -> 5983
$ws.Cells.Item(88, 9).Value = 5466.6665  # I88: 5650 -> 5466.6665
$ws.Cells.Item(88, 10).Value = 6241.1665  # J88: 6324.5 -> 6241.1665
$ws.Cells.Item(88, 11).Value = 5466.6665  # K88: 5650 -> 5466.6665
$ws.Cells.Item(88, 12).Value = 6241.1665  # L88: 6324.5 -> 6241.1665
$ws.Cells.Item(88, 13).Value = -5060.6665  # M88: -5244 -> -5060.6665
$ws.Cells.Item(88, 14).Value = -7053.1665  # N88: -7136.5 -> -7053.1665
$ws.Cells.Item(89, 8).Value = 5148.9443  # H89: 5042.684 -> 5148.9443
$ws.Cells.Item(89, 9).Value = 4697.125  # I89: 4523 -> 4697.125
$ws.Cells.Item(89, 11).Value = 23485.625  # K89: 22615 -> 23485.625
$ws.Cells.Item(89, 13).Value = -17869.625  # M89: -16999 -> -17869.625
$ws.Cells.Item(91, 8).Value = 5983  # H91: 6099.6665 -> 5983
$ws.Cells.Item(91, 9).Value = 5466.6665  # I91: 5650 -> 5466.6665
$ws.Cells.Item(91, 10).Value = 6241.1665  # J91: 6324.5 -> 6241.1665
$ws.Cells.Item(91, 11).Value = 5466.6665  # K91: 5650 -> 5466.6665
$ws.Cells.Item(91, 12).Value = 6241.1665  # L91: 6324.5 -> 6241.1665
$ws.Cells.Item(91, 13).Value = -4062.6665  # M91: -4246 -> -4062.6665
$ws.Cells.Item(91, 14).Value = -9049.1665  # N91: -9132.5 -> -9049.1665
$ws.Cells.Item(135, 8).Value = 3968.9678  # H135: 4077.9666 -> 3968.9678
$ws.Cells.Item(135, 9).Value = 3529.4075  # I135: 3638.2693 -> 3529.4075
$ws.Cells.Item(135, 11).Value = 31764.6675  # K135: 32744.4237 -> 31764.6675
$ws.Cells.Item(135, 13).Value = -29229.6675  # M135: -30209.4237 -> -29229.6675
$ws.Cells.Item(137, 8).Value = 6670006.5  # H137: 6899837.5 -> 6670006.5
$ws.Cells.Item(137, 9).Value = 11113210  # I137: 10002004 -> 11113210
$ws.Cells.Item(137, 10).Value = 5201.4165  # J137: 6135.222 -> 5201.4165
$ws.Cells.Item(137, 11).Value = 33339630  # K137: 30006012 -> 33339630
$ws.Cells.Item(137, 12).Value = 15604.2495  # L137: 18405.666 -> 15604.2495
$ws.Cells.Item(137, 13).Value = -33337080  # M137: -30003462 -> -33337080
$ws.Cells.Item(137, 14).Value = -20704.2495  # N137: -23505.666 -> -20704.2495
$ws.Cells.Item(138, 8).Value = 3560.2363  # H138: 3542.8147 -> 3560.2363
$ws.Cells.Item(138, 9).Value = 2520.4194  # I138: 2520.3872 -> 2520.4194
$ws.Cells.Item(138, 10).Value = 4903.3335  # J138: 4920.8696 -> 4903.3335
$ws.Cells.Item(138, 11).Value = 7561.2582  # K138: 7561.1616 -> 7561.2582
$ws.Cells.Item(138, 12).Value = 14710.0005  # L138: 14762.6088 -> 14710.0005
$ws.Cells.Item(138, 13).Value = -2421.2582  # M138: -2421.1616 -> -2421.2582
$ws.Cells.Item(138, 14).Value = -24990.0005  # N138: -25042.6088 -> -24990.0005
$ws.Cells.Item(141, 8).Value = 1714.5312  # H141: 1897.8 -> 1714.5312
$ws.Cells.Item(141, 9).Value = 1662.4445  # I141: 1884.3636 -> 1662.4445
$ws.Cells.Item(141, 10).Value = 1995.8  # J141: 1996.3334 -> 1995.8
$ws.Cells.Item(141, 11).Value = 4987.333500000001  # K141: 5653.0908 -> 4987.333500000001
$ws.Cells.Item(141, 12).Value = 5987.4  # L141: 5989.0002 -> 5987.4
$ws.Cells.Item(141, 13).Value = 192.6664999999994  # M141: -473.0907999999999 -> 192.6664999999994
$ws.Cells.Item(141, 14).Value = -16347.4  # N141: -16349.0002 -> -16347.4

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(5, 8).Value = 65.8125  # H5: 64.5 -> 65.8125
$ws.Cells.Item(5, 9).Value = 59.615383  # I5: 56.25 -> 59.615383
$ws.Cells.Item(5, 10).Value = 92.666664  # J5: 89.25 -> 92.666664
$ws.Cells.Item(5, 11).Value = 59.615383  # K5: 56.25 -> 59.615383
$ws.Cells.Item(5, 12).Value = 92.666664  # L5: 89.25 -> 92.666664
$ws.Cells.Item(5, 13).Value = 52.384617  # M5: 55.75 -> 52.384617
$ws.Cells.Item(5, 14).Value = -316.666664  # N5: -313.25 -> -316.666664
$ws.Cells.Item(32, 8).Value = 507070.7  # H32: 507072 -> 507070.7
$ws.Cells.Item(32, 10).Value = 2876.1667  # J32: 2892.8333 -> 2876.1667
$ws.Cells.Item(32, 12).Value = 2876.1667  # L32: 2892.8333 -> 2876.1667
$ws.Cells.Item(32, 14).Value = -3450.1667  # N32: -3466.8333 -> -3450.1667
$ws.Cells.Item(61, 8).Value = 31324512  # H61: 35364510 -> 31324512
$ws.Cells.Item(61, 9).Value = 30304370  # I61: 30637032 -> 30304370
$ws.Cells.Item(61, 10).Value = 31834582  # J61: 38201000 -> 31834582
$ws.Cells.Item(61, 11).Value = 30304370  # K61: 30637032 -> 30304370
$ws.Cells.Item(61, 12).Value = 31834582  # L61: 38201000 -> 31834582
$ws.Cells.Item(61, 13).Value = -30304158  # M61: -30636820 -> -30304158
$ws.Cells.Item(61, 14).Value = -31835006  # N61: -38201424 -> -31835006
$ws.Cells.Item(74, 8).Value = 1363063  # H74: 1246456.4 -> 1363063
$ws.Cells.Item(74, 9).Value = 2073704.6  # I74: 1893515.5 -> 2073704.6
$ws.Cells.Item(74, 10).Value = 6383.273  # J74: 6259.6665 -> 6383.273
$ws.Cells.Item(74, 11).Value = 2073704.6  # K74: 1893515.5 -> 2073704.6
$ws.Cells.Item(74, 12).Value = 6383.273  # L74: 6259.6665 -> 6383.273
$ws.Cells.Item(74, 13).Value = -2072830.6  # M74: -1892641.5 -> -2072830.6
$ws.Cells.Item(74, 14).Value = -8131.273  # N74: -8007.6665 -> -8131.273
$ws.Cells.Item(77, 8).Value = 1363063  # H77: 1246456.4 -> 1363063
$ws.Cells.Item(77, 9).Value = 2073704.6  # I77: 1893515.5 -> 2073704.6
$ws.Cells.Item(77, 10).Value = 6383.273  # J77: 6259.6665 -> 6383.273
$ws.Cells.Item(77, 11).Value = 10368523  # K77: 9467577.5 -> 10368523
$ws.Cells.Item(77, 12).Value = 31916.365  # L77: 31298.3325 -> 31916.365
$ws.Cells.Item(77, 13).Value = -10364155  # M77: -9463209.5 -> -10364155
$ws.Cells.Item(77, 14).Value = -40652.36500000001  # N77: -40034.3325 -> -40652.36500000001
$ws.Cells.Item(122, 8).Value = 1553.28  # H122: 1910.3684 -> 1553.28
$ws.Cells.Item(122, 9).Value = 1515.3636  # I122: 1870.4117 -> 1515.3636
$ws.Cells.Item(122, 10).Value = 1831.3334  # J122: 2250 -> 1831.3334
$ws.Cells.Item(122, 11).Value = 4546.0908  # K122: 5611.2351 -> 4546.0908
$ws.Cells.Item(122, 12).Value = 5494.0002  # L122: 6750 -> 5494.0002
$ws.Cells.Item(122, 13).Value = -2096.0908  # M122: -3161.2351 -> -2096.0908
$ws.Cells.Item(122, 14).Value = -10394.0002  # N122: -11650 -> -10394.0002
$ws.Cells.Item(132, 8).Value = 6079.653  # H132: 6281.2446 -> 6079.653
$ws.Cells.Item(132, 9).Value = 4413.9  # I132: 4542.9473 -> 4413.9
$ws.Cells.Item(132, 10).Value = 7228.448  # J132: 7551.5386 -> 7228.448
$ws.Cells.Item(132, 11).Value = 13241.7  # K132: 13628.8419 -> 13241.7
$ws.Cells.Item(132, 12).Value = 21685.344  # L132: 22654.6158 -> 21685.344
$ws.Cells.Item(132, 13).Value = -10711.7  # M132: -11098.8419 -> -10711.7
$ws.Cells.Item(132, 14).Value = -26745.344  # N132: -27714.6158 -> -26745.344
$ws.Cells.Item(133, 8).Value = 52954.4  # H133: 56904.2 -> 52954.4
$ws.Cells.Item(133, 10).Value = 52954.4  # J133: 56904.2 -> 52954.4
$ws.Cells.Item(133, 12).Value = 52954.4  # L133: 56904.2 -> 52954.4
$ws.Cells.Item(133, 14).Value = -58014.4  # N133: -61964.2 -> -58014.4
$ws.Cells.Item(136, 8).Value = 31324512  # H136: 35364510 -> 31324512
$ws.Cells.Item(136, 9).Value = 30304370  # I136: 30637032 -> 30304370
$ws.Cells.Item(136, 10).Value = 31834582  # J136: 38201000 -> 31834582
$ws.Cells.Item(136, 11).Value = 90913110  # K136: 91911096 -> 90913110
$ws.Cells.Item(136, 12).Value = 95503746  # L136: 114603000 -> 95503746
$ws.Cells.Item(136, 13).Value = -90910560  # M136: -91908546 -> -90910560
$ws.Cells.Item(136, 14).Value = -95508846  # N136: -114608100 -> -95508846
$ws.Cells.Item(137, 8).Value = 85259.29  # H137: 85306.336 -> 85259.29
$ws.Cells.Item(137, 10).Value = 85259.29  # J137: 85306.336 -> 85259.29
$ws.Cells.Item(137, 12).Value = 85259.29  # L137: 85306.336 -> 85259.29
$ws.Cells.Item(137, 14).Value = -95459.29  # N137: -95506.336 -> -95459.29
$ws.Cells.Item(139, 8).Value = 106506.5  # H139: 126631.664 -> 106506.5
$ws.Cells.Item(139, 10).Value = 106506.5  # J139: 126631.664 -> 106506.5
$ws.Cells.Item(139, 12).Value = 106506.5  # L139: 126631.664 -> 106506.5
$ws.Cells.Item(139, 14).Value = -116786.5  # N139: -136911.664 -> -116786.5

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(4, 8).Value = 65.8125  # H4: 64.5 -> 65.8125
$ws.Cells.Item(4, 9).Value = 59.615383  # I4: 56.25 -> 59.615383
$ws.Cells.Item(4, 10).Value = 92.666664  # J4: 89.25 -> 92.666664
$ws.Cells.Item(4, 11).Value = 59.615383  # K4: 56.25 -> 59.615383
$ws.Cells.Item(4, 12).Value = 92.666664  # L4: 89.25 -> 92.666664
$ws.Cells.Item(4, 13).Value = 55.384617  # M4: 58.75 -> 55.384617
$ws.Cells.Item(4, 14).Value = -322.666664  # N4: -319.25 -> -322.666664
$ws.Cells.Item(20, 8).Value = 652851.44  # H20: 761193.44 -> 652851.44
$ws.Cells.Item(20, 9).Value = 908692.8  # I20: 1135166.1 -> 908692.8
$ws.Cells.Item(20, 11).Value = 908692.8  # K20: 1135166.1 -> 908692.8
$ws.Cells.Item(20, 13).Value = -908445.8  # M20: -1134919.1 -> -908445.8
$ws.Cells.Item(22, 8).Value = 750.1429  # H22: 748 -> 750.1429
$ws.Cells.Item(22, 9).Value = 750.1429  # I22: 833.5 -> 750.1429
$ws.Cells.Item(22, 10).Value = 0  # J22: 235 -> 0
$ws.Cells.Item(22, 11).Value = 750.1429  # K22: 833.5 -> 750.1429
$ws.Cells.Item(22, 12).Value = 0  # L22: 235 -> 0
$ws.Cells.Item(22, 13).Value = -577.1429  # M22: -660.5 -> -577.1429
$ws.Cells.Item(22, 14).Value = $null  # N22: -581 -> None
$ws.Cells.Item(58, 8).Value = 71199.8  # H58: 61118 -> 71199.8
$ws.Cells.Item(58, 9).Value = 0  # I58: 5709 -> 0
$ws.Cells.Item(58, 10).Value = 71199.8  # J58: 72199.8 -> 71199.8
$ws.Cells.Item(58, 11).Value = 0  # K58: 5709 -> 0
$ws.Cells.Item(58, 12).Value = 71199.8  # L58: 72199.8 -> 71199.8
$ws.Cells.Item(58, 13).Value = $null  # M58: -5415 -> None
$ws.Cells.Item(58, 14).Value = -71787.8  # N58: -72787.8 -> -71787.8
$ws.Cells.Item(59, 8).Value = 76405.8  # H59: 93493 -> 76405.8
$ws.Cells.Item(59, 10).Value = 76405.8  # J59: 93493 -> 76405.8
$ws.Cells.Item(59, 12).Value = 76405.8  # L59: 93493 -> 76405.8
$ws.Cells.Item(59, 14).Value = -78099.8  # N59: -95187 -> -78099.8
$ws.Cells.Item(86, 8).Value = 4577.4287  # H86: 4504.1387 -> 4577.4287
$ws.Cells.Item(86, 9).Value = 4289.4814  # I86: 4389.077 -> 4289.4814
$ws.Cells.Item(86, 10).Value = 5549.25  # J86: 4803.3 -> 5549.25
$ws.Cells.Item(86, 11).Value = 4289.4814  # K86: 4389.077 -> 4289.4814
$ws.Cells.Item(86, 12).Value = 5549.25  # L86: 4803.3 -> 5549.25
$ws.Cells.Item(86, 13).Value = -3166.4814  # M86: -3266.077 -> -3166.4814
$ws.Cells.Item(86, 14).Value = -7795.25  # N86: -7049.3 -> -7795.25
$ws.Cells.Item(89, 8).Value = 4577.4287  # H89: 4504.1387 -> 4577.4287
$ws.Cells.Item(89, 9).Value = 4289.4814  # I89: 4389.077 -> 4289.4814
$ws.Cells.Item(89, 10).Value = 5549.25  # J89: 4803.3 -> 5549.25
$ws.Cells.Item(89, 11).Value = 21447.407  # K89: 21945.385 -> 21447.407
$ws.Cells.Item(89, 12).Value = 27746.25  # L89: 24016.5 -> 27746.25
$ws.Cells.Item(89, 13).Value = -15831.407  # M89: -16329.385 -> -15831.407
$ws.Cells.Item(89, 14).Value = -38978.25  # N89: -35248.5 -> -38978.25
$ws.Cells.Item(99, 8).Value = 14238.5  # H99: 15598.333 -> 14238.5
$ws.Cells.Item(99, 9).Value = 17355  # I99: 19914.166 -> 17355
$ws.Cells.Item(99, 11).Value = 17355  # K99: 19914.166 -> 17355
$ws.Cells.Item(99, 13).Value = -15857  # M99: -18416.166 -> -15857
$ws.Cells.Item(105, 8).Value = 4072.682  # H105: 4116.6113 -> 4072.682
$ws.Cells.Item(105, 9).Value = 3740.2666  # I105: 3864.5715 -> 3740.2666
$ws.Cells.Item(105, 10).Value = 4785  # J105: 4998.75 -> 4785
$ws.Cells.Item(105, 11).Value = 3740.2666  # K105: 3864.5715 -> 3740.2666
$ws.Cells.Item(105, 12).Value = 4785  # L105: 4998.75 -> 4785
$ws.Cells.Item(105, 13).Value = -1993.2666  # M105: -2117.5715 -> -1993.2666
$ws.Cells.Item(105, 14).Value = -8279  # N105: -8492.75 -> -8279
$ws.Cells.Item(114, 8).Value = 20310.5  # H114: 30000 -> 20310.5
$ws.Cells.Item(114, 9).Value = 10621  # I114: 0 -> 10621
$ws.Cells.Item(114, 11).Value = 10621  # K114: 0 -> 10621
$ws.Cells.Item(114, 13).Value = -6282  # M114: None -> -6282
$ws.Cells.Item(132, 8).Value = 95780  # H132: 95578 -> 95780
$ws.Cells.Item(132, 10).Value = 95780  # J132: 95578 -> 95780
$ws.Cells.Item(132, 12).Value = 95780  # L132: 95578 -> 95780
$ws.Cells.Item(132, 14).Value = -105900  # N132: -105698 -> -105900
$ws.Cells.Item(134, 8).Value = 4632259.5  # H134: 5170649 -> 4632259.5
$ws.Cells.Item(134, 9).Value = 4168623  # I134: 4387983 -> 4168623
$ws.Cells.Item(134, 10).Value = 6950443  # J134: 11118911 -> 6950443
$ws.Cells.Item(134, 11).Value = 12505869  # K134: 13163949 -> 12505869
$ws.Cells.Item(134, 12).Value = 20851329  # L134: 33356733 -> 20851329
$ws.Cells.Item(134, 13).Value = -12503334  # M134: -13161414 -> -12503334
$ws.Cells.Item(134, 14).Value = -20856399  # N134: -33361803 -> -20856399

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 60617.65  # H16: 64344.312 -> 60617.65
$ws.Cells.Item(16, 9).Value = 1695.7693  # I16: 1754.5 -> 1695.7693
$ws.Cells.Item(16, 11).Value = 1695.7693  # K16: 1754.5 -> 1695.7693
$ws.Cells.Item(16, 13).Value = -1408.7693  # M16: -1467.5 -> -1408.7693
$ws.Cells.Item(22, 8).Value = 1439.2727  # H22: 1567.2 -> 1439.2727
$ws.Cells.Item(22, 9).Value = 851.2857  # I22: 966.5 -> 851.2857
$ws.Cells.Item(22, 11).Value = 851.2857  # K22: 966.5 -> 851.2857
$ws.Cells.Item(22, 13).Value = -501.2857  # M22: -616.5 -> -501.2857
$ws.Cells.Item(31, 8).Value = 752505.94  # H31: 795857.75 -> 752505.94
$ws.Cells.Item(31, 9).Value = 1424221.1  # I31: 1529644 -> 1424221.1
$ws.Cells.Item(31, 10).Value = 3285.2307  # J31: 3368.64 -> 3285.2307
$ws.Cells.Item(31, 11).Value = 1424221.1  # K31: 1529644 -> 1424221.1
$ws.Cells.Item(31, 12).Value = 3285.2307  # L31: 3368.64 -> 3285.2307
$ws.Cells.Item(31, 13).Value = -1423926.1  # M31: -1529349 -> -1423926.1
$ws.Cells.Item(31, 14).Value = -3875.2307  # N31: -3958.64 -> -3875.2307
$ws.Cells.Item(34, 8).Value = 752505.94  # H34: 795857.75 -> 752505.94
$ws.Cells.Item(34, 9).Value = 1424221.1  # I34: 1529644 -> 1424221.1
$ws.Cells.Item(34, 10).Value = 3285.2307  # J34: 3368.64 -> 3285.2307
$ws.Cells.Item(34, 11).Value = 1424221.1  # K34: 1529644 -> 1424221.1
$ws.Cells.Item(34, 12).Value = 3285.2307  # L34: 3368.64 -> 3285.2307
$ws.Cells.Item(34, 13).Value = -1424019.1  # M34: -1529442 -> -1424019.1
$ws.Cells.Item(34, 14).Value = -3689.2307  # N34: -3772.64 -> -3689.2307
$ws.Cells.Item(41, 8).Value = 19799  # H41: 20899 -> 19799
$ws.Cells.Item(58, 8).Value = 26272480  # H58: 19107676 -> 26272480
$ws.Cells.Item(58, 9).Value = 33338316  # I58: 23813382 -> 33338316
$ws.Cells.Item(58, 10).Value = 14496087  # J58: 10872690 -> 14496087
$ws.Cells.Item(58, 11).Value = 33338316  # K58: 23813382 -> 33338316
$ws.Cells.Item(58, 12).Value = 14496087  # L58: 10872690 -> 14496087
$ws.Cells.Item(58, 13).Value = -33338113  # M58: -23813179 -> -33338113
$ws.Cells.Item(58, 14).Value = -14496493  # N58: -10873096 -> -14496493
$ws.Cells.Item(86, 8).Value = 9238.741  # H86: 9618.84 -> 9238.741
$ws.Cells.Item(86, 9).Value = 5710  # I86: 5856.4287 -> 5710
$ws.Cells.Item(86, 10).Value = 10473.8  # J86: 11082 -> 10473.8
$ws.Cells.Item(86, 11).Value = 5710  # K86: 5856.4287 -> 5710
$ws.Cells.Item(86, 12).Value = 10473.8  # L86: 11082 -> 10473.8
$ws.Cells.Item(86, 13).Value = -4587  # M86: -4733.4287 -> -4587
$ws.Cells.Item(86, 14).Value = -12719.8  # N86: -13328 -> -12719.8
$ws.Cells.Item(89, 8).Value = 9238.741  # H89: 9618.84 -> 9238.741
$ws.Cells.Item(89, 9).Value = 5710  # I89: 5856.4287 -> 5710
$ws.Cells.Item(89, 10).Value = 10473.8  # J89: 11082 -> 10473.8
$ws.Cells.Item(89, 11).Value = 28550  # K89: 29282.1435 -> 28550
$ws.Cells.Item(89, 12).Value = 52369  # L89: 55410 -> 52369
$ws.Cells.Item(89, 13).Value = -22934  # M89: -23666.1435 -> -22934
$ws.Cells.Item(89, 14).Value = -63601  # N89: -66642 -> -63601
$ws.Cells.Item(105, 8).Value = 8402.212  # H105: 8637.4375 -> 8402.212
$ws.Cells.Item(105, 9).Value = 5820.826  # I105: 6045.636 -> 5820.826
$ws.Cells.Item(105, 11).Value = 5820.826  # K105: 6045.636 -> 5820.826
$ws.Cells.Item(105, 13).Value = -4073.826  # M105: -4298.636 -> -4073.826
$ws.Cells.Item(113, 8).Value = 60617.65  # H113: 64344.312 -> 60617.65
$ws.Cells.Item(113, 9).Value = 1695.7693  # I113: 1754.5 -> 1695.7693
$ws.Cells.Item(113, 11).Value = 1695.7693  # K113: 1754.5 -> 1695.7693
$ws.Cells.Item(113, 13).Value = 474.2307000000001  # M113: 415.5 -> 474.2307000000001
$ws.Cells.Item(122, 8).Value = 15006.1875  # H122: 18246.615 -> 15006.1875
$ws.Cells.Item(122, 9).Value = 2476.8462  # I122: 2930.6 -> 2476.8462
$ws.Cells.Item(122, 11).Value = 7430.5386  # K122: 8791.8 -> 7430.5386
$ws.Cells.Item(122, 13).Value = -4980.5386  # M122: -6341.799999999999 -> -4980.5386
$ws.Cells.Item(132, 8).Value = 3843.577  # H132: 4084.6667 -> 3843.577
$ws.Cells.Item(132, 9).Value = 3651.3809  # I132: 3935.6843 -> 3651.3809
$ws.Cells.Item(132, 11).Value = 10954.1427  # K132: 11807.0529 -> 10954.1427
$ws.Cells.Item(132, 13).Value = -8424.1427  # M132: -9277.052899999999 -> -8424.1427
$ws.Cells.Item(134, 8).Value = 4304.548  # H134: 4385.244 -> 4304.548
$ws.Cells.Item(134, 9).Value = 2046.5294  # I134: 2112.1875 -> 2046.5294
$ws.Cells.Item(134, 11).Value = 6139.5882  # K134: 6336.5625 -> 6139.5882
$ws.Cells.Item(134, 13).Value = -3604.5882  # M134: -3801.5625 -> -3604.5882
$ws.Cells.Item(136, 8).Value = 26272480  # H136: 19107676 -> 26272480
$ws.Cells.Item(136, 9).Value = 33338316  # I136: 23813382 -> 33338316
$ws.Cells.Item(136, 10).Value = 14496087  # J136: 10872690 -> 14496087
$ws.Cells.Item(136, 11).Value = 100014948  # K136: 71440146 -> 100014948
$ws.Cells.Item(136, 12).Value = 43488261  # L136: 32618070 -> 43488261
$ws.Cells.Item(136, 13).Value = -100012398  # M136: -71437596 -> -100012398
$ws.Cells.Item(136, 14).Value = -43493361  # N136: -32623170 -> -43493361

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 3575295  # H4: 3575295.5 -> 3575295
$ws.Cells.Item(4, 9).Value = 4085888.8  # I4: 4085889.2 -> 4085888.8
$ws.Cells.Item(4, 11).Value = 12257666.4  # K4: 12257667.6 -> 12257666.4
$ws.Cells.Item(4, 13).Value = -12257554.4  # M4: -12257555.6 -> -12257554.4
$ws.Cells.Item(6, 8).Value = 1190.1  # H6: 1687.1428 -> 1190.1
$ws.Cells.Item(6, 9).Value = 71.2  # I6: 107 -> 71.2
$ws.Cells.Item(6, 10).Value = 2309  # J6: 2872.25 -> 2309
$ws.Cells.Item(6, 11).Value = 213.6  # K6: 321 -> 213.6
$ws.Cells.Item(6, 12).Value = 6927  # L6: 8616.75 -> 6927
$ws.Cells.Item(6, 13).Value = -100.6  # M6: -208 -> -100.6
$ws.Cells.Item(6, 14).Value = -7153  # N6: -8842.75 -> -7153
$ws.Cells.Item(10, 8).Value = 1257.25  # H10: 1673.3334 -> 1257.25
$ws.Cells.Item(10, 9).Value = 9.666667  # I10: 10 -> 9.666667
$ws.Cells.Item(10, 11).Value = 29.000001  # K10: 30 -> 29.000001
$ws.Cells.Item(10, 13).Value = 109.999999  # M10: 109 -> 109.999999
$ws.Cells.Item(11, 8).Value = 83333470  # H11: 90909240 -> 83333470
$ws.Cells.Item(11, 9).Value = 121.666664  # I11: 124.75 -> 121.666664
$ws.Cells.Item(11, 11).Value = 364.999992  # K11: 374.25 -> 364.999992
$ws.Cells.Item(11, 13).Value = -224.999992  # M11: -234.25 -> -224.999992
$ws.Cells.Item(17, 8).Value = 666.8889  # H17: 653.5 -> 666.8889
$ws.Cells.Item(17, 9).Value = 100  # I17: 298.7143 -> 100
$ws.Cells.Item(17, 10).Value = 1120.4  # J17: 1150.2 -> 1120.4
$ws.Cells.Item(17, 11).Value = 300  # K17: 896.1428999999999 -> 300
$ws.Cells.Item(17, 12).Value = 3361.2  # L17: 3450.6 -> 3361.2
$ws.Cells.Item(17, 13).Value = -131  # M17: -727.1428999999999 -> -131
$ws.Cells.Item(17, 14).Value = -3699.2  # N17: -3788.6 -> -3699.2
$ws.Cells.Item(68, 8).Value = 3847.0908  # H68: 3812.513 -> 3847.0908
$ws.Cells.Item(68, 10).Value = 4031  # J68: 3991.5342 -> 4031
$ws.Cells.Item(68, 12).Value = 12093  # L68: 11974.6026 -> 12093
$ws.Cells.Item(68, 14).Value = -13715  # N68: -13596.6026 -> -13715
$ws.Cells.Item(71, 8).Value = 3847.0908  # H71: 3812.513 -> 3847.0908
$ws.Cells.Item(71, 10).Value = 4031  # J71: 3991.5342 -> 4031
$ws.Cells.Item(71, 12).Value = 36279  # L71: 35923.8078 -> 36279
$ws.Cells.Item(71, 14).Value = -44391  # N71: -44035.8078 -> -44391
$ws.Cells.Item(92, 8).Value = 1000  # H92: 3000 -> 1000
$ws.Cells.Item(92, 9).Value = 1000  # I92: 3000 -> 1000
$ws.Cells.Item(92, 10).Value = 1000  # J92: 0 -> 1000
$ws.Cells.Item(92, 11).Value = 3000  # K92: 9000 -> 3000
$ws.Cells.Item(92, 12).Value = 3000  # L92: 0 -> 3000
$ws.Cells.Item(92, 13).Value = -1752  # M92: -7752 -> -1752
$ws.Cells.Item(92, 14).Value = -5496  # N92: None -> -5496
$ws.Cells.Item(107, 8).Value = 3385.4092  # H107: 3644.4 -> 3385.4092
$ws.Cells.Item(107, 9).Value = 756.6667  # I107: 719.8571 -> 756.6667
$ws.Cells.Item(107, 10).Value = 4371.1875  # J107: 5219.154 -> 4371.1875
$ws.Cells.Item(107, 11).Value = 2270.0001  # K107: 2159.5713 -> 2270.0001
$ws.Cells.Item(107, 12).Value = 13113.5625  # L107: 15657.462 -> 13113.5625
$ws.Cells.Item(107, 13).Value = -350.0001000000002  # M107: -239.5712999999996 -> -350.0001000000002
$ws.Cells.Item(107, 14).Value = -16953.5625  # N107: -19497.462 -> -16953.5625
$ws.Cells.Item(113, 8).Value = 469  # H113: 483.69565 -> 469
$ws.Cells.Item(113, 10).Value = 508.29413  # J113: 536.06665 -> 508.29413
$ws.Cells.Item(113, 12).Value = 1524.88239  # L113: 1608.19995 -> 1524.88239
$ws.Cells.Item(113, 14).Value = -5864.88239  # N113: -5948.19995 -> -5864.88239
$ws.Cells.Item(131, 8).Value = 7026.4614  # H131: 7758.8184 -> 7026.4614
$ws.Cells.Item(131, 10).Value = 7797.2  # J131: 8996.875 -> 7797.2
$ws.Cells.Item(131, 12).Value = 23391.6  # L131: 26990.625 -> 23391.6
$ws.Cells.Item(131, 14).Value = -33471.6  # N131: -37070.625 -> -33471.6
$ws.Cells.Item(132, 8).Value = 1155.9166  # H132: 972.4167 -> 1155.9166
$ws.Cells.Item(132, 9).Value = 634.8  # I132: 481.57144 -> 634.8
$ws.Cells.Item(132, 10).Value = 1528.1428  # J132: 1659.6 -> 1528.1428
$ws.Cells.Item(132, 11).Value = 5713.2  # K132: 4334.14296 -> 5713.2
$ws.Cells.Item(132, 12).Value = 13753.2852  # L132: 14936.4 -> 13753.2852
$ws.Cells.Item(132, 13).Value = -3183.2  # M132: -1804.14296 -> -3183.2
$ws.Cells.Item(132, 14).Value = -18813.2852  # N132: -19996.4 -> -18813.2852
$ws.Cells.Item(138, 8).Value = 60074.668  # H138: 45305.875 -> 60074.668
$ws.Cells.Item(138, 9).Value = 105233  # I138: 63539.6 -> 105233
$ws.Cells.Item(138, 11).Value = 315699  # K138: 190618.8 -> 315699
$ws.Cells.Item(138, 13).Value = -310559  # M138: -185478.8 -> -310559

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 3580  # H80: 4140 -> 3580
$ws.Cells.Item(80, 9).Value = 1975  # I80: 2233.3333 -> 1975
$ws.Cells.Item(80, 10).Value = 10000  # J80: 7000 -> 10000
$ws.Cells.Item(80, 11).Value = 1975  # K80: 2233.3333 -> 1975
$ws.Cells.Item(80, 12).Value = 10000  # L80: 7000 -> 10000
$ws.Cells.Item(80, 13).Value = -977  # M80: -1235.3333 -> -977
$ws.Cells.Item(80, 14).Value = -11996  # N80: -8996 -> -11996
$ws.Cells.Item(83, 8).Value = 3580  # H83: 4140 -> 3580
$ws.Cells.Item(83, 9).Value = 1975  # I83: 2233.3333 -> 1975
$ws.Cells.Item(83, 10).Value = 10000  # J83: 7000 -> 10000
$ws.Cells.Item(83, 11).Value = 9875  # K83: 11166.6665 -> 9875
$ws.Cells.Item(83, 12).Value = 50000  # L83: 35000 -> 50000
$ws.Cells.Item(83, 13).Value = -4883  # M83: -6174.666499999999 -> -4883
$ws.Cells.Item(83, 14).Value = -59984  # N83: -44984 -> -59984
$ws.Cells.Item(102, 8).Value = 1489.9117  # H102: 1623.7333 -> 1489.9117
$ws.Cells.Item(102, 9).Value = 1568.2667  # I102: 1734.7307 -> 1568.2667
$ws.Cells.Item(102, 11).Value = 1568.2667  # K102: 1734.7307 -> 1568.2667
$ws.Cells.Item(102, 13).Value = 53.7333000000001  # M102: -112.7307000000001 -> 53.7333000000001
$ws.Cells.Item(122, 8).Value = 33255.305  # H122: 36524.832 -> 33255.305
$ws.Cells.Item(122, 9).Value = 53176.75  # I122: 59046.945 -> 53176.75
$ws.Cells.Item(122, 10).Value = 2606.923  # J122: 2741.6667 -> 2606.923
$ws.Cells.Item(122, 11).Value = 159530.25  # K122: 177140.835 -> 159530.25
$ws.Cells.Item(122, 12).Value = 7820.768999999999  # L122: 8225.000100000001 -> 7820.768999999999
$ws.Cells.Item(122, 13).Value = -157080.25  # M122: -174690.835 -> -157080.25
$ws.Cells.Item(122, 14).Value = -12720.769  # N122: -13125.0001 -> -12720.769
$ws.Cells.Item(132, 8).Value = 10903.932  # H132: 11629.22 -> 10903.932
$ws.Cells.Item(132, 9).Value = 7650.387  # I132: 8363.821 -> 7650.387
$ws.Cells.Item(132, 11).Value = 22951.161  # K132: 25091.463 -> 22951.161
$ws.Cells.Item(132, 13).Value = -20421.161  # M132: -22561.463 -> -20421.161

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 2775.85  # H22: 2905.85 -> 2775.85
$ws.Cells.Item(22, 9).Value = 2177.3333  # I22: 2324.5 -> 2177.3333
$ws.Cells.Item(22, 10).Value = 3265.5454  # J22: 3293.4167 -> 3265.5454
$ws.Cells.Item(22, 11).Value = 2177.3333  # K22: 2324.5 -> 2177.3333
$ws.Cells.Item(22, 12).Value = 3265.5454  # L22: 3293.4167 -> 3265.5454
$ws.Cells.Item(22, 13).Value = -1882.3333  # M22: -2029.5 -> -1882.3333
$ws.Cells.Item(22, 14).Value = -3855.5454  # N22: -3883.4167 -> -3855.5454
$ws.Cells.Item(27, 8).Value = 2775.85  # H27: 2905.85 -> 2775.85
$ws.Cells.Item(27, 9).Value = 2177.3333  # I27: 2324.5 -> 2177.3333
$ws.Cells.Item(27, 10).Value = 3265.5454  # J27: 3293.4167 -> 3265.5454
$ws.Cells.Item(27, 11).Value = 2177.3333  # K27: 2324.5 -> 2177.3333
$ws.Cells.Item(27, 12).Value = 3265.5454  # L27: 3293.4167 -> 3265.5454
$ws.Cells.Item(27, 13).Value = -2070.3333  # M27: -2217.5 -> -2070.3333
$ws.Cells.Item(27, 14).Value = -3479.5454  # N27: -3507.4167 -> -3479.5454
$ws.Cells.Item(46, 8).Value = 3640.5356  # H46: 3724.4443 -> 3640.5356
$ws.Cells.Item(46, 9).Value = 3099  # I46: 3083.4285 -> 3099
$ws.Cells.Item(46, 10).Value = 3788.2273  # J46: 3948.8 -> 3788.2273
$ws.Cells.Item(46, 11).Value = 3099  # K46: 3083.4285 -> 3099
$ws.Cells.Item(46, 12).Value = 3788.2273  # L46: 3948.8 -> 3788.2273
$ws.Cells.Item(46, 13).Value = -2911  # M46: -2895.4285 -> -2911
$ws.Cells.Item(46, 14).Value = -4164.2273  # N46: -4324.8 -> -4164.2273
$ws.Cells.Item(93, 8).Value = 2092.1  # H93: 2097.1 -> 2092.1
$ws.Cells.Item(93, 9).Value = 2101.889  # I93: 2107.4443 -> 2101.889
$ws.Cells.Item(93, 11).Value = 2101.889  # K93: 2107.4443 -> 2101.889
$ws.Cells.Item(93, 13).Value = -853.8890000000001  # M93: -859.4443000000001 -> -853.8890000000001
$ws.Cells.Item(100, 8).Value = 2219  # H100: 2992.3333 -> 2219
$ws.Cells.Item(100, 9).Value = 2219  # I100: 2488.5 -> 2219
$ws.Cells.Item(100, 10).Value = 0  # J100: 4000 -> 0
$ws.Cells.Item(100, 11).Value = 2219  # K100: 2488.5 -> 2219
$ws.Cells.Item(100, 12).Value = 0  # L100: 4000 -> 0
$ws.Cells.Item(100, 13).Value = -1678  # M100: -1947.5 -> -1678
$ws.Cells.Item(100, 14).Value = $null  # N100: -5082 -> None
$ws.Cells.Item(132, 8).Value = 2567744.8  # H132: 2567802 -> 2567744.8
$ws.Cells.Item(132, 9).Value = 3707386.8  # I132: 4170460 -> 3707386.8
$ws.Cells.Item(132, 10).Value = 3550  # J132: 3549 -> 3550
$ws.Cells.Item(132, 11).Value = 11122160.4  # K132: 12511380 -> 11122160.4
$ws.Cells.Item(132, 12).Value = 10650  # L132: 10647 -> 10650
$ws.Cells.Item(132, 13).Value = -11119630.4  # M132: -12508850 -> -11119630.4
$ws.Cells.Item(132, 14).Value = -15710  # N132: -15707 -> -15710
$ws.Cells.Item(133, 8).Value = 120500  # H133: 110333 -> 120500
$ws.Cells.Item(133, 10).Value = 120500  # J133: 110333 -> 120500
$ws.Cells.Item(133, 12).Value = 120500  # L133: 110333 -> 120500
$ws.Cells.Item(133, 14).Value = -125560  # N133: -115393 -> -125560
$ws.Cells.Item(136, 8).Value = 10649948  # H136: 12779672 -> 10649948
$ws.Cells.Item(136, 9).Value = 8335293.5  # I136: 10418784 -> 8335293.5
$ws.Cells.Item(136, 11).Value = 25005880.5  # K136: 31256352 -> 25005880.5
$ws.Cells.Item(136, 13).Value = -25003330.5  # M136: -31253802 -> -25003330.5
$ws.Cells.Item(137, 8).Value = 73949.4  # H137: 73955.4 -> 73949.4
$ws.Cells.Item(137, 10).Value = 73949.4  # J137: 73955.4 -> 73949.4
$ws.Cells.Item(137, 12).Value = 73949.4  # L137: 73955.4 -> 73949.4
$ws.Cells.Item(137, 14).Value = -84149.4  # N137: -84155.4 -> -84149.4
$ws.Cells.Item(139, 8).Value = 90620  # H139: 90630 -> 90620
$ws.Cells.Item(139, 10).Value = 90620  # J139: 90630 -> 90620
$ws.Cells.Item(139, 12).Value = 90620  # L139: 90630 -> 90620
$ws.Cells.Item(139, 14).Value = -100900  # N139: -100910 -> -100900

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(46, 8).Value = 62143  # H46: 61476.332 -> 62143
$ws.Cells.Item(46, 10).Value = 62143  # J46: 61476.332 -> 62143
$ws.Cells.Item(46, 12).Value = 62143  # L46: 61476.332 -> 62143
$ws.Cells.Item(46, 14).Value = -62605  # N46: -61938.332 -> -62605
$ws.Cells.Item(62, 8).Value = 11164.632  # H62: 11453.947 -> 11164.632
$ws.Cells.Item(62, 9).Value = 5718.75  # I62: 6176.8335 -> 5718.75
$ws.Cells.Item(62, 11).Value = 5718.75  # K62: 6176.8335 -> 5718.75
$ws.Cells.Item(62, 13).Value = -5094.75  # M62: -5552.8335 -> -5094.75
$ws.Cells.Item(65, 8).Value = 11164.632  # H65: 11453.947 -> 11164.632
$ws.Cells.Item(65, 9).Value = 5718.75  # I65: 6176.8335 -> 5718.75
$ws.Cells.Item(65, 11).Value = 28593.75  # K65: 30884.1675 -> 28593.75
$ws.Cells.Item(65, 13).Value = -25473.75  # M65: -27764.1675 -> -25473.75
$ws.Cells.Item(74, 8).Value = 20600  # H74: 20599.5 -> 20600
$ws.Cells.Item(74, 10).Value = 20600  # J74: 20599.5 -> 20600
$ws.Cells.Item(74, 12).Value = 20600  # L74: 20599.5 -> 20600
$ws.Cells.Item(74, 14).Value = -22472  # N74: -22471.5 -> -22472
$ws.Cells.Item(77, 8).Value = 20600  # H77: 20599.5 -> 20600
$ws.Cells.Item(77, 10).Value = 20600  # J77: 20599.5 -> 20600
$ws.Cells.Item(77, 12).Value = 61800  # L77: 61798.5 -> 61800
$ws.Cells.Item(77, 14).Value = -71160  # N77: -71158.5 -> -71160
$ws.Cells.Item(122, 8).Value = 38191.06  # H122: 39815.453 -> 38191.06
$ws.Cells.Item(122, 9).Value = 3909.862  # I122: 4600 -> 3909.862
$ws.Cells.Item(122, 11).Value = 11729.586  # K122: 13800 -> 11729.586
$ws.Cells.Item(122, 13).Value = -9279.586  # M122: -11350 -> -9279.586
$ws.Cells.Item(126, 8).Value = 2879.1785  # H126: 2948.889 -> 2879.1785
$ws.Cells.Item(126, 9).Value = 2008.3043  # I126: 2054.2727 -> 2008.3043
$ws.Cells.Item(126, 11).Value = 6024.9129  # K126: 6162.8181 -> 6024.9129
$ws.Cells.Item(126, 13).Value = -3554.9129  # M126: -3692.8181 -> -3554.9129
$ws.Cells.Item(132, 8).Value = 8337034.5  # H132: 9263105 -> 8337034.5
$ws.Cells.Item(132, 9).Value = 10420193  # I132: 11908449 -> 10420193
$ws.Cells.Item(132, 11).Value = 31260579  # K132: 35725347 -> 31260579
$ws.Cells.Item(132, 13).Value = -31258049  # M132: -35722817 -> -31258049
$ws.Cells.Item(134, 8).Value = 62143  # H134: 61476.332 -> 62143
$ws.Cells.Item(134, 10).Value = 62143  # J134: 61476.332 -> 62143
$ws.Cells.Item(134, 12).Value = 186429  # L134: 184428.996 -> 186429
$ws.Cells.Item(134, 14).Value = -191499  # N134: -189498.996 -> -191499
